$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.185.33"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").Value = "2.057.51"
$ws.Range("E3").Value = "  +0.78%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.72"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.93%  "
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.27"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.60%  "
$ws.Range("E9").Value = "  +0.52%  "
$ws.Range("E10").Value = "  +0.14%  "
$ws.Range("E11").Value = "  +0.66%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.23"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.71%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.914"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +13.35%  "
$ws.Range("D14").Value = "2.360.17"
$ws.Range("E14").Value = "  +0.82%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.75"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.28%  "
$ws.Range("D16").Value = "2.064.82"
$ws.Range("E16").Value = "  +1.17%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.78"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +11.75%  "
$ws.Range("D18").Value = "37.194.86"
$ws.Range("E18").Value = "  +0.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "74.91"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.22%  "
$ws.Range("D20").Value = "0.0₃0899"
$ws.Range("E20").Value = "  -0.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.48"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "238.18"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.59%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.49"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +4.55%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.64"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +4.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.19"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -3.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "170.58"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.23"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.31%  "
$ws.Range("E29").Value = "  -0.30%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.15"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +8.99%  "
$ws.Range("E31").Value = "  +3.30%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0625"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.85%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.67"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +4.22%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0883"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.24%  "
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("E36").Value = "  +2.37%  "
$ws.Range("E37").Value = "  +1.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.34"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.23"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +12.90%  "
$ws.Range("E40").Value = "  -8.63%  "
$ws.Range("E41").Value = "  +7.64%  "
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.56"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.40%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.16"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +2.57%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "96.79"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.42"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.03%  "
$ws.Range("D47").Value = "1.278.83"
$ws.Range("E47").Value = "  -0.30%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.86"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.91%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.85"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.36%  "
$ws.Range("D50").Value = "2.247.56"
$ws.Range("E50").Value = "  +0.97%  "
$ws.Range("E51").Value = "  +9.62%  "
